$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize capitalization of connector words ("de" -> "De", "y" -> "Y") in specific municipality/state names
$ws.Range("B2").Value = "Comitán De Domínguez"
$ws.Range("A12").Value = "Ciudad De México"
$ws.Range("A17").Value = "Estado De México"
$ws.Range("B18").Value = "Naucalpan De Juárez"
$ws.Range("B19").Value = "Tlalnepantla De Baz"
$ws.Range("B23").Value = "Tlapa De Comonfort"
$ws.Range("B25").Value = "Jacala De Ledezma"
$ws.Range("B33").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B46").Value = "Martínez De La Torre"
$ws.Range("B48").Value = "Sayula De Alemán"

# Remove trailing metadata/footer rows (53-57)
$ws.Range("A53:A57").EntireRow.Delete()
